$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Result (%)" column with recalculated figures based on 46 total test cases
$ws.Range("E5").Value = "(46/46)*100 = 100"
$ws.Range("E6").Value = "(0/46)*100 = 0"
$ws.Range("E7").Value = "(44/46)*100 = 95.65"
$ws.Range("E8").Value = "(2/46)*100 = 4.34"
$ws.Range("E9").Value = "(0/46)*100 = 0"

# Update the active selection to match the author's final cursor position
$ws.Range("E14").Select()
